$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param([int]$Row, [object[]]$Values)
    $data = New-Object 'object[,]' 1,10
    for ($i = 0; $i -lt 10; $i++) {
        $data[0,$i] = $Values[$i]
    }
    $startCell = $ws.Cells.Item($Row, 4)
    $endCell = $ws.Cells.Item($Row, 13)
    $ws.Range($startCell, $endCell).Value2 = $data
}

# Step 1: extend formatting from column K into new columns L:M for each data block,
# so the new cells inherit the correct style (matching a fresh column insert).
$ws.Range("K7:K35").Copy() | Out-Null
$ws.Range("L7:M35").PasteSpecial(-4122) | Out-Null
$ws.Range("K38:K77").Copy() | Out-Null
$ws.Range("L38:M77").PasteSpecial(-4122) | Out-Null
$ws.Range("K80:K102").Copy() | Out-Null
$ws.Range("L80:M102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Step 2: write the refreshed financial data (two additional quarters prepended,
# plus revised historical figures) across D:M for every data row.

Set-RowData 7 @(43465,43373,43281,43190,43100,43008,42916,42825,42735,42643)
Set-RowData 8 @(253700,258900,263800,261900,257200,235900,243400,249300,236500,232600)
Set-RowData 9 @(185200,185800,189300,183300,180000,165700,165100,164100,157900,155900)
Set-RowData 10 @(68500,73100,74500,78600,77200,70200,78300,85200,78600,76700)
Set-RowData 11 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
Set-RowData 12 @(35300,19100,21000,23300,33900,19700,21400,19500,25000,19700)
Set-RowData 13 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 14 @(1900,17300,6200,900,800,0,0,0,100,0)
Set-RowData 15 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 16 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
Set-RowData 17 @(232800,254700,248200,241200,236200,219700,218200,214400,210100,205200)
Set-RowData 18 @(20900,4200,15600,20700,21000,16200,25200,34900,26400,27400)
Set-RowData 19 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
Set-RowData 20 @(-600,400,5400,-3500,-1200,-7700,-13000,-1100,6900,-500)
Set-RowData 21 @(32200,17500,33900,30100,32100,19900,23100,43900,43400,37100)
Set-RowData 22 @(1300,1200,1200,1200,1300,1300,1300,1100,1000,700)
Set-RowData 23 @(19000,3300,19700,16000,18600,7200,10900,32600,32400,26200)
Set-RowData 24 @(6400,3700,3100,3000,3700,600,2400,7200,6300,6000)
Set-RowData 25 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 26 @(12600,-400,16700,13000,14900,6600,8500,25400,26000,20200)
Set-RowData 27 @(12600,-400,16700,13000,14900,6600,8500,25400,26000,20200)
Set-RowData 28 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 29 @("NA","NA","NA","NA",-20100,"NA","NA","NA","NA","NA")
Set-RowData 30 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 31 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 32 @(600,-400,-5400,3500,1200,7700,13000,1100,-6900,500)
Set-RowData 33 @(12600,-400,16700,13000,-5200,6600,8500,25400,26000,20200)
Set-RowData 34 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 35 @(12600,-400,16700,13000,-5200,6600,8500,25400,26000,20200)
Set-RowData 38 @(43465,43373,43281,43190,43100,43008,42916,42825,42735,42643)
Set-RowData 39 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
Set-RowData 40 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
Set-RowData 41 @(39600,47200,65400,70500,103200,147600,164200,133900,177200,132800)
Set-RowData 42 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 43 @(166900,183500,200000,196600,185100,182200,181300,185300,170100,179000)
Set-RowData 44 @(112500,108600,118400,120600,121400,113400,106600,108500,105100,107200)
Set-RowData 45 @(124100,115300,62800,58000,51400,38600,42500,45500,32000,46000)
Set-RowData 46 @(443100,454500,446600,445800,461100,481800,494500,473100,484400,465000)
Set-RowData 47 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 48 @(171400,189000,203900,205200,200300,190800,187400,179800,172100,167300)
Set-RowData 49 @(111700,115900,142400,151200,153000,106800,107900,107300,109300,115900)
Set-RowData 50 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 51 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 52 @(76900,85600,97100,99200,69100,83100,76800,72500,77300,61300)
Set-RowData 53 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 54 @(803000,845000,890000,901400,883400,862600,866500,832700,843000,809600)
Set-RowData 55 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
Set-RowData 56 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
Set-RowData 57 @(93100,89000,95000,88100,89600,80300,83700,86200,84500,88800)
Set-RowData 58 @(3400,3400,3400,3500,3500,3400,3400,2100,2100,900)
Set-RowData 59 @(78900,86000,73200,72900,78300,65500,65900,64900,107000,112400)
Set-RowData 60 @(175400,178400,171700,164500,171300,149300,153100,153200,193600,202100)
Set-RowData 61 @(136500,98000,109500,105700,141200,142400,159900,161000,169400,140700)
Set-RowData 62 @(11500,18300,20400,20800,17000,20100,20600,21200,19600,21700)
Set-RowData 63 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 64 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 65 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 66 @(323300,294700,301600,291000,329500,311800,333500,335400,382600,364400)
Set-RowData 67 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
Set-RowData 68 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 69 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 70 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 71 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 72 @(364000,351300,351700,335000,293600,298900,292300,283800,256900,230900)
Set-RowData 73 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 74 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 75 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 76 @(479700,550300,588400,610400,553900,550800,533000,497300,460400,445100)
Set-RowData 77 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 80 @(43465,43373,43281,43190,43100,43008,42916,42825,42735,42643)
Set-RowData 81 @(12600,-400,16700,13000,-5200,6600,8500,25400,26000,20200)
Set-RowData 82 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
Set-RowData 83 @(11900,12900,12900,12900,12300,11500,11000,10200,10000,10200)
Set-RowData 84 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 85 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 86 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 87 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 88 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 89 @(47900,38000,27000,5600,34500,9300,27500,-21300,37400,19800)
Set-RowData 90 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
Set-RowData 91 @(-9700,-9700,-13800,-8400,-13600,-11400,-12200,-13600,-15600,-19900)
Set-RowData 92 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 93 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 94 @(-9600,-9700,-13100,-8400,-78500,-11400,-12200,-15600,-20000,-19800)
Set-RowData 95 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
Set-RowData 96 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 97 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 98 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 99 @(0,0,0,0,0,0,0,0,0,0)
Set-RowData 100 @(-45100,-46300,-12500,-35400,-1700,-22400,1000,-8500,38300,-300)
Set-RowData 101 @(-700,-200,-6500,5500,1300,8000,14000,2100,-11300,1100)
Set-RowData 102 @(-7500,-18200,-5100,-32700,-44500,-16600,30300,-43300,44400,800)

